# Commit: "Doing Updates for Financials"
# Adds the two newest quarters (period ending 2018-12-31 and 2018-09-30) to every
# table on the AZN sheet, by inserting two new data columns right after column C
# and filling in the new figures, shifting the previously-existing quarters over.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D; this shifts the existing D:K data to F:M.
$ws.Range("D:E").Insert()

# Populate the two new columns: D = quarter ending 2018-12-31, E = quarter ending 2018-09-30.
$ws.Cells.Item(7, 4).Value2 = 43465
$ws.Cells.Item(7, 5).Value2 = 43373
$ws.Cells.Item(8, 4).Value = 6417000
$ws.Cells.Item(8, 5).Value = 5340000
$ws.Cells.Item(9, 4).Value = 1282000
$ws.Cells.Item(9, 5).Value = 1131000
$ws.Cells.Item(10, 4).Value = 5135000
$ws.Cells.Item(10, 5).Value = 4209000
$ws.Cells.Item(12, 4).Value = 2005000
$ws.Cells.Item(12, 5).Value = 1242000
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 46000
$ws.Cells.Item(14, 5).Value = 191000
$ws.Cells.Item(15, 4).Value = 883000
$ws.Cells.Item(15, 5).Value = 373000
$ws.Cells.Item(17, 4).Value = 5340000
$ws.Cells.Item(17, 5).Value = 4489000
$ws.Cells.Item(18, 4).Value = 1077000
$ws.Cells.Item(18, 5).Value = 851000
$ws.Cells.Item(20, 4).Value = -347000
$ws.Cells.Item(20, 5).Value = -374000
$ws.Cells.Item(21, 4).Value = 2392000
$ws.Cells.Item(21, 5).Value = 1175000
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 4).Value = 730000
$ws.Cells.Item(23, 5).Value = 477000
$ws.Cells.Item(24, 4).Value = -279000
$ws.Cells.Item(24, 5).Value = 71000
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 1009000
$ws.Cells.Item(26, 5).Value = 406000
$ws.Cells.Item(27, 4).Value = 1034000
$ws.Cells.Item(27, 5).Value = 431000
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = "NA"
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = 347000
$ws.Cells.Item(32, 5).Value = 374000
$ws.Cells.Item(33, 4).Value = 1034000
$ws.Cells.Item(33, 5).Value = 431000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 1034000
$ws.Cells.Item(35, 5).Value = 431000
$ws.Cells.Item(38, 4).Value2 = 43465
$ws.Cells.Item(38, 5).Value2 = 43373
$ws.Cells.Item(41, 4).Value = 4831000
$ws.Cells.Item(41, 5).Value = 3420000
$ws.Cells.Item(42, 4).Value = 849000
$ws.Cells.Item(42, 5).Value = 808000
$ws.Cells.Item(43, 4).Value = 5781000
$ws.Cells.Item(43, 5).Value = 5819000
$ws.Cells.Item(44, 4).Value = 2890000
$ws.Cells.Item(44, 5).Value = 3027000
$ws.Cells.Item(45, 4).Value = 1240000
$ws.Cells.Item(45, 5).Value = 34000
$ws.Cells.Item(46, 4).Value = 15591000
$ws.Cells.Item(46, 5).Value = 13108000
$ws.Cells.Item(47, 4).Value = 1437000
$ws.Cells.Item(47, 5).Value = 1942000
$ws.Cells.Item(48, 4).Value = 7421000
$ws.Cells.Item(48, 5).Value = 7591000
$ws.Cells.Item(49, 4).Value = 33666000
$ws.Cells.Item(49, 5).Value = 36147000
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 2536000
$ws.Cells.Item(52, 5).Value = 2655000
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 60651000
$ws.Cells.Item(54, 5).Value = 61443000
$ws.Cells.Item(57, 4).Value = 12841000
$ws.Cells.Item(57, 5).Value = 10992000
$ws.Cells.Item(58, 4).Value = 1754000
$ws.Cells.Item(58, 5).Value = 2491000
$ws.Cells.Item(59, 4).Value = 1697000
$ws.Cells.Item(59, 5).Value = 1765000
$ws.Cells.Item(60, 4).Value = 16292000
$ws.Cells.Item(60, 5).Value = 15248000
$ws.Cells.Item(61, 4).Value = 17359000
$ws.Cells.Item(61, 5).Value = 18422000
$ws.Cells.Item(62, 4).Value = 12956000
$ws.Cells.Item(62, 5).Value = 14236000
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 48183000
$ws.Cells.Item(66, 5).Value = 49507000
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 7724000
$ws.Cells.Item(72, 5).Value = 7202000
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 12468000
$ws.Cells.Item(76, 5).Value = 11936000
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value2 = 43465
$ws.Cells.Item(80, 5).Value2 = 43373
$ws.Cells.Item(81, 4).Value = 1034000
$ws.Cells.Item(81, 5).Value = 431000
$ws.Cells.Item(83, 4).Value = 1662000
$ws.Cells.Item(83, 5).Value = 698000
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 2224000
$ws.Cells.Item(89, 5).Value = 469000
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = 927000
$ws.Cells.Item(94, 5).Value = -141000
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = -1121000
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -1732000
$ws.Cells.Item(100, 5).Value = 169000
$ws.Cells.Item(101, 4).Value = -10000
$ws.Cells.Item(101, 5).Value = -1000
$ws.Cells.Item(102, 4).Value = 1409000
$ws.Cells.Item(102, 5).Value = 496000

# Row 91 "Capital Expenditures" carries revised figures across the whole row (D:M),
# not merely the two new quarters, so set all ten values explicitly.
$ws.Cells.Item(91, 4).Value = -315000
$ws.Cells.Item(91, 5).Value = -242000
$ws.Cells.Item(91, 6).Value = -273000
$ws.Cells.Item(91, 7).Value = -213000
$ws.Cells.Item(91, 8).Value = -477000
$ws.Cells.Item(91, 9).Value = -300000
$ws.Cells.Item(91, 10).Value = -263000
$ws.Cells.Item(91, 11).Value = -385000
$ws.Cells.Item(91, 12).Value = -641000
$ws.Cells.Item(91, 13).Value = -366000
